$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.889.83'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '1.755.38'
$ws.Range('E3').Value = '  -4.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9982'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5106'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.31'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2791'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06227'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -10.18%  '
$ws.Range('D11').Value = '1.749.04'
$ws.Range('E11').Value = '  -4.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.87'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.14%  '
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6137'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -16.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.533'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -9.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -12.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9989'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9990'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = '25.874.61'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007029'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -11.08%  '
$ws.Range('E21').Value = '  -15.24%  '
$ws.Range('D22').Value = '1.970.88'
$ws.Range('E22').Value = '  -5.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.092'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.272'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -12.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.240'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -10.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.474'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -13.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.835'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -15.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '103.87'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08217'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.723'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.535'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -12.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04533'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9979'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -10.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9980'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -11.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6152'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -15.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.700'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -12.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01560'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.916'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -16.65%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '103.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9988'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3891'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -17.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7494'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -16.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.948'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -15.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05404'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1117'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.035'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -18.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -12.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -11.80%  '
